$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
